# Added comments and cleaned the output
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shuffle / update the seating order in column B (rows 2-25).
$ws.Range("B2").Value  = "Vera"
$ws.Range("B3").Value  = "David"
$ws.Range("B4").Value  = "An"
$ws.Range("B5").Value  = "Kevin J"
$ws.Range("B6").Value  = "Jessica"
$ws.Range("B7").Value  = "Andrii"
$ws.Range("B8").Value  = "Miriam"
$ws.Range("B9").Value  = "Imad"
$ws.Range("B10").Value = "Patrick"
$ws.Range("B11").Value = "Fatemeh"
$ws.Range("B12").Value = "Mohamad"
$ws.Range("B13").Value = "Patrycja"
$ws.Range("B14").Value = "Therese"
$ws.Range("B15").Value = "Miro"
$ws.Range("B16").Value = "Frank"
$ws.Range("B17").Value = "Olha"
$ws.Range("B18").Value = "Beatriz"
$ws.Range("B19").Value = "Dhanya"
$ws.Range("B20").Value = "Manel"
$ws.Range("B21").Value = "Celina"
$ws.Range("B22").Value = "Aleksander"
$ws.Range("B23").Value = "Oscar"
$ws.Range("B24").Value = "Edoardo"
$ws.Range("B25").Value = "Kevin P"
